# Continue to change structure
# Adds 3 new condition rows (16-18) to the condition.csv sheet, matching the
# new "commerce invest money", "trade condition" and "unblock trade item"
# conditions, and resizes columns A/D to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: canCommerceInvestMoneyEnough ---------------------------------
$ws.Range("A16").Value = "canCommerceInvestMoneyEnough"
$ws.Range("B16").Value = "商业投资钱是否够"
$ws.Range("C16").Value = "city"
$ws.Range("D16").Value = "commerceInvestMoney"
$ws.Range("E16").Value = "<="
$ws.Range("F16").Value = "guild"
$ws.Range("G16").Value = "money"

# --- Row 17: tradeCondition -------------------------------------------------
$ws.Range("A17").Value = "tradeCondition"
$ws.Range("B17").Value = "签约了"
$ws.Range("C17").Value = "city"
$ws.Range("D17").Value = "percentage"
$ws.Range("E17").Value = "'>"
$ws.Range("F17").Value = "'number"
$ws.Range("G17").Value = 0

# --- Row 18: canUnblockTradeItem -------------------------------------------
$ws.Range("A18").Value = "canUnblockTradeItem"
$ws.Range("B18").Value = "可以解锁商品"
$ws.Range("C18").Value = "cacheString"
$ws.Range("D18").Value = "tradeUnblockItem"
$ws.Range("E18").Value = "!="
$ws.Range("F18").Value = "number"
$ws.Range("G18").Value = 0

# --- Resize columns A and D to fit the new, wider content ------------------
$ws.Columns.Item(1).ColumnWidth = 28.666666666666664
$ws.Columns.Item(4).ColumnWidth = 19.5

# --- Selection ends on the last edited cell, like the authored commit ------
$ws.Range("D18").Select() | Out-Null
